$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.96000000000062"
$ws.Range("H2").Value = [double]"0.3935938753219833"
$ws.Range("I2").Value = [double]"0.3935938753219833"
$ws.Range("L2").Value = [double]"3.770348528762438"
$ws.Range("N2").Value = [double]"0.3019407041416662"
$ws.Range("O2").Value = [double]"0.3019407041416662"
$ws.Range("P2").Value = [double]"-1.74847398938531"
$ws.Range("R2").Value = [double]"0.2639831432516635"
$ws.Range("S2").Value = [double]"0.2639831432516635"
$ws.Range("T2").Value = [double]"11.27589802921457"
$ws.Range("V2").Value = [double]"9.276035546346151e-07"
$ws.Range("W2").Value = [double]"9.276035546346151e-07"
$ws.Range("X2").Value = [double]"7.224104104104274"
$ws.Range("Y2").Value = [double]"-5.638958958959096"
$ws.Range("Z2").Value = [double]"20.08716716716764"
$ws.Range("M2").Value = "[-3.5016359251390274, 11.042332982663904]"
$ws.Range("Q2").Value = "[-4.861764006456276, 1.3648160276856558]"
$ws.Range("U2").Value = "[7.278042130895699, 15.273753927533438]"

# Row 3
$ws.Range("F3").Value = [double]"25.96000000000062"
$ws.Range("H3").Value = [double]"0.2684950861458723"
$ws.Range("I3").Value = [double]"0.2684950861458723"
$ws.Range("L3").Value = [double]"4.979876216103064"
$ws.Range("N3").Value = [double]"0.227699749898097"
$ws.Range("O3").Value = [double]"0.227699749898097"
$ws.Range("P3").Value = [double]"-2.050368778919465"
$ws.Range("R3").Value = [double]"0.169260067517701"
$ws.Range("S3").Value = [double]"0.169260067517701"
$ws.Range("T3").Value = [double]"14.98985300309567"
$ws.Range("V3").Value = [double]"1.128324877974762e-08"
$ws.Range("W3").Value = [double]"1.128324877974762e-08"
$ws.Range("X3").Value = [double]"8.471431431431633"
$ws.Range("Y3").Value = [double]"-3.74198198198207"
$ws.Range("Z3").Value = [double]"20.68484484484534"
$ws.Range("M3").Value = "[-3.2213232141504458, 13.181075646356575]"
$ws.Range("Q3").Value = "[-5.006421926441392, 0.9056843686024623]"
$ws.Range("U3").Value = "[10.658368020681323, 19.321337985510027]"

# Row 4
$ws.Range("F4").Value = [double]"25.96000000000062"
$ws.Range("H4").Value = [double]"0.2401514149513491"
$ws.Range("I4").Value = [double]"0.2401514149513491"
$ws.Range("L4").Value = [double]"5.748857212806155"
$ws.Range("N4").Value = [double]"0.2058644930736093"
$ws.Range("O4").Value = [double]"0.2058644930736093"
$ws.Range("P4").Value = [double]"-1.861684535460618"
$ws.Range("R4").Value = [double]"0.229131560696914"
$ws.Range("S4").Value = [double]"0.229131560696914"
$ws.Range("T4").Value = [double]"16.21853220907132"
$ws.Range("V4").Value = [double]"2.6334890712576e-08"
$ws.Range("W4").Value = [double]"2.6334890712576e-08"
$ws.Range("X4").Value = [double]"7.691851851852036"
$ws.Range("Y4").Value = [double]"-5.015295295295418"
$ws.Range("Z4").Value = [double]"20.39899899899949"
$ws.Range("M4").Value = "[-3.2719460617826206, 14.76966048739493]"
$ws.Range("Q4").Value = "[-4.937237703839816, 1.2138686329185795]"
$ws.Range("U4").Value = "[11.359195610029817, 21.077868808112825]"

# Row 5
$ws.Range("F5").Value = [double]"25.96000000000062"
$ws.Range("H5").Value = [double]"0.1697377072380168"
$ws.Range("I5").Value = [double]"0.1697377072380168"
$ws.Range("L5").Value = [double]"5.12460561578918"
$ws.Range("N5").Value = [double]"0.1401292563762306"
$ws.Range("O5").Value = [double]"0.1401292563762306"
$ws.Range("P5").Value = [double]"2.861711025792505"
$ws.Range("R5").Value = [double]"0.03356442982683161"
$ws.Range("S5").Value = [double]"0.03356442982683161"
$ws.Range("T5").Value = [double]"11.48688782690343"
$ws.Range("V5").Value = [double]"1.713925688040518e-07"
$ws.Range("W5").Value = [double]"1.713925688040518e-07"
$ws.Range("X5").Value = [double]"14.13637637637671"
$ws.Range("Y5").Value = [double]"3.27423423423431"
$ws.Range("Z5").Value = [double]"24.99851851851911"
$ws.Range("M5").Value = "[-1.7480542040126625, 11.997265435591023]"
$ws.Range("Q5").Value = "[0.2327105669325773, 5.490711484652432]"
$ws.Range("U5").Value = "[7.740412235516738, 15.23336341829013]"

# Row 6
$ws.Range("B6").Value = 0
$ws.Range("F6").Value = [double]"25.96000000000062"
$ws.Range("H6").Value = [double]"0.8970361765767813"
$ws.Range("I6").Value = [double]"0.8970361765767813"
$ws.Range("L6").Value = [double]"1.414768995320921"
$ws.Range("N6").Value = [double]"0.7255348908259318"
$ws.Range("O6").Value = [double]"0.7255348908259318"
$ws.Range("P6").Value = [double]"-2.717053105807389"
$ws.Range("R6").Value = [double]"0.0874292414409743"
$ws.Range("S6").Value = [double]"0.0874292414409743"
$ws.Range("T6").Value = [double]"10.74056753740831"
$ws.Range("V6").Value = [double]"4.431934654602898e-06"
$ws.Range("W6").Value = [double]"4.431934654602898e-06"
$ws.Range("X6").Value = [double]"11.22594594594621"
$ws.Range("Y6").Value = [double]"-1.715075075075116"
$ws.Range("Z6").Value = [double]"24.16696696696754"
$ws.Range("M6").Value = "[-6.651250522936611, 9.480788513578453]"
$ws.Range("Q6").Value = "[-5.84921154722424, 0.4151053356094625]"
$ws.Range("U6").Value = "[6.595063609332897, 14.886071465483731]"

# Row 7
$ws.Range("F7").Value = [double]"25.96000000000062"
$ws.Range("H7").Value = [double]"0.06319193831835257"
$ws.Range("I7").Value = [double]"0.06319193831835257"
$ws.Range("L7").Value = [double]"7.381539983089645"
$ws.Range("N7").Value = [double]"0.09575657936296711"
$ws.Range("O7").Value = [double]"0.09575657936296711"
$ws.Range("P7").Value = [double]"2.522079387566581"
$ws.Range("R7").Value = [double]"0.06991036117224536"
$ws.Range("S7").Value = [double]"0.06991036117224536"
$ws.Range("T7").Value = [double]"13.96982528483336"
$ws.Range("V7").Value = [double]"8.066661227879024e-08"
$ws.Range("W7").Value = [double]"8.066661227879024e-08"
$ws.Range("X7").Value = [double]"15.53961961961999"
$ws.Range("Y7").Value = [double]"4.235715715715813"
$ws.Range("Z7").Value = [double]"26.84352352352416"
$ws.Range("M7").Value = "[-1.3563345356089798, 16.11941450178827]"
$ws.Range("Q7").Value = "[-0.21384214258669232, 5.258000917719855]"
$ws.Range("U7").Value = "[9.570313631981225, 18.3693369376855]"

# Row 8
$ws.Range("F8").Value = [double]"25.96000000000062"
$ws.Range("H8").Value = [double]"0.3813034416379262"
$ws.Range("I8").Value = [double]"0.3813034416379262"
$ws.Range("L8").Value = [double]"4.864631993034623"
$ws.Range("N8").Value = [double]"0.3313589890667346"
$ws.Range("O8").Value = [double]"0.3313589890667346"
$ws.Range("P8").Value = [double]"2.03150035457358"
$ws.Range("R8").Value = [double]"0.1980760697920929"
$ws.Range("S8").Value = [double]"0.1980760697920929"
$ws.Range("T8").Value = [double]"17.25066281182394"
$ws.Range("V8").Value = [double]"2.374917529301968e-08"
$ws.Range("W8").Value = [double]"2.374917529301968e-08"
$ws.Range("X8").Value = [double]"17.56652652652695"
$ws.Range("Y8").Value = [double]"4.625505505505616"
$ws.Range("Z8").Value = [double]"30.50754754754828"
$ws.Range("M8").Value = "[-5.113016362889784, 14.84228034895903]"
$ws.Range("Q8").Value = "[-1.1006580868432714, 5.163658795990431]"
$ws.Range("U8").Value = "[12.105198073887259, 22.396127549760628]"

# Row 9
$ws.Range("F9").Value = [double]"22.80000000000013"
$ws.Range("H9").Value = [double]"0.3881366143283921"
$ws.Range("I9").Value = [double]"0.3881366143283921"
$ws.Range("L9").Value = [double]"4.516639622578596"
$ws.Range("N9").Value = [double]"0.2775448449179285"
$ws.Range("O9").Value = [double]"0.2775448449179285"
$ws.Range("P9").Value = [double]"1.088079137279347"
$ws.Range("R9").Value = [double]"0.4868589405512931"
$ws.Range("S9").Value = [double]"0.4868589405512931"
$ws.Range("T9").Value = [double]"12.33064928317409"
$ws.Range("V9").Value = [double]"1.983504313018969e-06"
$ws.Range("W9").Value = [double]"1.983504313018969e-06"
$ws.Range("X9").Value = [double]"18.85165165165176"
$ws.Range("Y9").Value = [double]"7.50870870870876"
$ws.Range("Z9").Value = [double]"30.19459459459476"
$ws.Range("M9").Value = "[-3.7596726205692246, 12.792951865726417]"
$ws.Range("Q9").Value = "[-2.0377898293555416, 4.213948103914235]"
$ws.Range("U9").Value = "[7.779382248951293, 16.88191631739688]"

# Row 10
$ws.Range("F10").Value = [double]"22.80000000000013"
$ws.Range("H10").Value = [double]"0.04306966566079762"
$ws.Range("I10").Value = [double]"0.04306966566079762"
$ws.Range("L10").Value = [double]"7.901610431658655"
$ws.Range("N10").Value = [double]"0.04682988632101659"
$ws.Range("O10").Value = [double]"0.04682988632101659"
$ws.Range("P10").Value = [double]"1.314500229429963"
$ws.Range("R10").Value = [double]"0.06437051367951252"
$ws.Range("S10").Value = [double]"0.06437051367951252"
$ws.Range("T10").Value = [double]"11.27070660331649"
$ws.Range("V10").Value = [double]"2.775333504123623e-06"
$ws.Range("W10").Value = [double]"2.775333504123623e-06"
$ws.Range("X10").Value = [double]"18.03003003003013"
$ws.Range("Y10").Value = [double]"12.96336336336343"
$ws.Range("Z10").Value = [double]"23.09669669669683"
$ws.Range("M10").Value = "[0.11566331706314337, 15.687557546254167]"
$ws.Range("Q10").Value = "[-0.08176317216550189, 2.7107636310254275]"
$ws.Range("U10").Value = "[7.03347648344381, 15.507936723189161]"

# Row 11
$ws.Range("B11").Value = 1
$ws.Range("F11").Value = [double]"22.80000000000013"
$ws.Range("H11").Value = [double]"0.008208494257470722"
$ws.Range("I11").Value = [double]"0.008208494257470722"
$ws.Range("L11").Value = [double]"9.246163865069693"
$ws.Range("N11").Value = [double]"0.01044266991410803"
$ws.Range("O11").Value = [double]"0.01044266991410803"
$ws.Range("P11").Value = [double]"1.389973926813502"
$ws.Range("R11").Value = [double]"0.007866266447666748"
$ws.Range("S11").Value = [double]"0.007866266447666748"
$ws.Range("T11").Value = [double]"10.50413512942792"
$ws.Range("V11").Value = [double]"2.578468242475651e-06"
$ws.Range("W11").Value = [double]"2.578468242475651e-06"
$ws.Range("X11").Value = [double]"17.75615615615625"
$ws.Range("Y11").Value = [double]"14.10450450450458"
$ws.Range("Z11").Value = [double]"21.40780780780793"
$ws.Range("M11").Value = "[2.2784679291617103, 16.213859800977676]"
$ws.Range("Q11").Value = "[0.3836579616996545, 2.396289891927349]"
$ws.Range("U11").Value = "[6.571100901615916, 14.437169357239926]"

# Row 12
$ws.Range("F12").Value = [double]"22.80000000000013"
$ws.Range("H12").Value = [double]"0.2466470110482092"
$ws.Range("I12").Value = [double]"0.2466470110482092"
$ws.Range("L12").Value = [double]"5.636924193320492"
$ws.Range("N12").Value = [double]"0.2050452364864388"
$ws.Range("O12").Value = [double]"0.2050452364864388"
$ws.Range("P12").Value = [double]"2.295658295415965"
$ws.Range("R12").Value = [double]"0.1444770930476968"
$ws.Range("S12").Value = [double]"0.1444770930476968"
$ws.Range("T12").Value = [double]"14.62245958382096"
$ws.Range("V12").Value = [double]"1.478380409203339e-07"
$ws.Range("W12").Value = [double]"1.478380409203339e-07"
$ws.Range("X12").Value = [double]"14.46966966966975"
$ws.Range("Y12").Value = [double]"3.17237237237239"
$ws.Range("Z12").Value = [double]"25.76696696696711"
$ws.Range("M12").Value = "[-3.191976550669062, 14.465824937310046]"
$ws.Range("Q12").Value = "[-0.8176317216550011, 5.408948312486931]"
$ws.Range("U12").Value = "[9.886429247267102, 19.358489920374808]"

# Row 13
$ws.Range("F13").Value = [double]"22.80000000000013"
$ws.Range("H13").Value = [double]"0.07285559864644298"
$ws.Range("I13").Value = [double]"0.07285559864644298"
$ws.Range("L13").Value = [double]"8.642986287617653"
$ws.Range("N13").Value = [double]"0.09248884319255968"
$ws.Range("O13").Value = [double]"0.09248884319255968"
$ws.Range("P13").Value = [double]"2.572395185822273"
$ws.Range("R13").Value = [double]"0.02016700166773577"
$ws.Range("S13").Value = [double]"0.02016700166773577"
$ws.Range("T13").Value = [double]"15.27942591698821"
$ws.Range("V13").Value = [double]"5.975409631542306e-07"
$ws.Range("W13").Value = [double]"5.975409631542306e-07"
$ws.Range("X13").Value = [double]"13.46546546546554"
$ws.Range("Y13").Value = [double]"5.660060060060093"
$ws.Range("Z13").Value = [double]"21.27087087087099"
$ws.Range("M13").Value = "[-1.483799704575155, 18.76977227981046]"
$ws.Range("Q13").Value = "[0.421394810391424, 4.723395561253122]"
$ws.Range("U13").Value = "[9.982531425251274, 20.576320408725138]"

# Row 14
$ws.Range("F14").Value = [double]"22.80000000000013"
$ws.Range("H14").Value = [double]"0.570971844431662"
$ws.Range("I14").Value = [double]"0.570971844431662"
$ws.Range("L14").Value = [double]"2.854000785463495"
$ws.Range("N14").Value = [double]"0.4276907341810725"
$ws.Range("O14").Value = [double]"0.4276907341810725"
$ws.Range("P14").Value = [double]"3.062974218815274"
$ws.Range("R14").Value = [double]"0.05319546890258819"
$ws.Range("S14").Value = [double]"0.05319546890258819"
$ws.Range("T14").Value = [double]"9.829627611501099"
$ws.Range("V14").Value = [double]"7.293181489620437e-06"
$ws.Range("W14").Value = [double]"7.293181489620437e-06"
$ws.Range("X14").Value = [double]"11.68528528528535"
$ws.Range("Y14").Value = [double]"0.4108108108108155"
$ws.Range("Z14").Value = [double]"22.95975975975988"
$ws.Range("M14").Value = "[-4.327819849990039, 10.035821420917028]"
$ws.Range("Q14").Value = "[-0.044026323473730145, 6.169974761104278]"
$ws.Range("U14").Value = "[5.924426067352431, 13.734829155649766]"
